# Fix Training Data Issue (#48)
# Data was taken from 1 day off due to way NBA stats were shown.
# Column BF ("Date") held values formatted as "5-3-2007-08" (a mangled
# mash-up of game-day and season strings); correct this to the proper
# ISO-like date string "2008-05-03" for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "5-3-2007-08"
$newValue = "2008-05-03"
$dateCol = 58  # Column BF

# Row 1 holds the "Date" header; data rows follow it through the bottom
# of the used range.
$lastRow = $ws.UsedRange.Rows.Count

# A leading apostrophe forces Excel to store the corrected value as
# literal text instead of auto-converting the ISO-like date string into a
# date serial number; re-applying the "Normal" style afterwards drops the
# transient quote-prefix formatting so the cell's style stays untouched.
for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, $dateCol)
    if ($cell.Value() -eq $oldValue) {
        $cell.Value = "'" + $newValue
        $cell.Style = "Normal"
    }
}
